$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item("Title 1")
$title1.TextFrame.TextRange.Text = "-"
$title1.TextFrame.TextRange.Text = "Slide 1"

$tb1 = $s1.Shapes.Item("TextBox 3")
$tb1.TextFrame.TextRange.Text = "-"
$tb1.TextFrame.TextRange.Text = "an image"

$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item("Title 1")
$title2.TextFrame.TextRange.Text = "-"
$title2.TextFrame.TextRange.Text = "Slide 2"

$tb2 = $s2.Shapes.Item("TextBox 3")
$tb2.TextFrame.TextRange.Text = "-"
$tb2.TextFrame.TextRange.Text = "an image"
